$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update res_bus vm_pu values for rows 2-25 (time steps 0-23), columns B-F and I-N
# (column G stays 1, column H has no data in these rows) - case with 380 kV done

# Row 2
$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.038937967464187
$ws.Range("D2").Value2 = 1.046162500309767
$ws.Range("E2").Value2 = 1.047375584475634
$ws.Range("F2").Value2 = 1.057893531192926
$ws.Range("I2").Value2 = 1.035640922432004
$ws.Range("J2").Value2 = 1.044032693845849
$ws.Range("K2").Value2 = 1.048928708558647
$ws.Range("L2").Value2 = 1.050138396932887
$ws.Range("M2").Value2 = 1.060627254425254
$ws.Range("N2").Value2 = 1.045515340023334

# Row 3
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.039937040499387
$ws.Range("D3").Value2 = 1.046938383800776
$ws.Range("E3").Value2 = 1.048256700037881
$ws.Range("F3").Value2 = 1.05884240883902
$ws.Range("I3").Value2 = 1.035796980430403
$ws.Range("J3").Value2 = 1.044676517946861
$ws.Range("K3").Value2 = 1.049516271783306
$ws.Range("L3").Value2 = 1.050831164369094
$ws.Range("M3").Value2 = 1.061389715462832
$ws.Range("N3").Value2 = 1.046160078428417

# Row 4
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.040584035743644
$ws.Range("D4").Value2 = 1.047440678817488
$ws.Range("E4").Value2 = 1.048827650595713
$ws.Range("F4").Value2 = 1.059457211112308
$ws.Range("I4").Value2 = 1.035896677779639
$ws.Range("J4").Value2 = 1.045093043218389
$ws.Range("K4").Value2 = 1.049896051391723
$ws.Range("L4").Value2 = 1.051279602522668
$ws.Range("M4").Value2 = 1.061883267275403
$ws.Range("N4").Value2 = 1.046577195213607

# Row 5
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.040856158237041
$ws.Range("D5").Value2 = 1.047651901359809
$ws.Range("E5").Value2 = 1.049067870872316
$ws.Range("F5").Value2 = 1.059715867805255
$ws.Range("I5").Value2 = 1.035938282980436
$ws.Range("J5").Value2 = 1.045268132333039
$ws.Range("K5").Value2 = 1.050055610852918
$ws.Range("L5").Value2 = 1.051468165949756
$ws.Range("M5").Value2 = 1.062090800476773
$ws.Range("N5").Value2 = 1.046752532974882

# Row 6
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.040901856149912
$ws.Range("D6").Value2 = 1.04768736990089
$ws.Range("E6").Value2 = 1.04910821617606
$ws.Range("F6").Value2 = 1.059759308730152
$ws.Range("I6").Value2 = 1.03594525062631
$ws.Range("J6").Value2 = 1.045297529478349
$ws.Range("K6").Value2 = 1.050082395702346
$ws.Range("L6").Value2 = 1.05149982889759
$ws.Range("M6").Value2 = 1.062125648763353
$ws.Range("N6").Value2 = 1.046781971867511

# Row 7
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.040587671364767
$ws.Range("D7").Value2 = 1.047443500958279
$ws.Range("E7").Value2 = 1.048830859675308
$ws.Range("F7").Value2 = 1.059460666535747
$ws.Range("I7").Value2 = 1.035897234919517
$ws.Range("J7").Value2 = 1.045095382839967
$ws.Range("K7").Value2 = 1.049898183826782
$ws.Range("L7").Value2 = 1.05128212196036
$ws.Range("M7").Value2 = 1.061886040172631
$ws.Range("N7").Value2 = 1.046579538157716

# Row 8
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.039275499054269
$ws.Range("D8").Value2 = 1.046424661944442
$ws.Range("E8").Value2 = 1.04767319318096
$ws.Range("F8").Value2 = 1.058214039310426
$ws.Range("I8").Value2 = 1.035693928232352
$ws.Range("J8").Value2 = 1.044250291588846
$ws.Range("K8").Value2 = 1.04912736323708
$ws.Range("L8").Value2 = 1.050372484872463
$ws.Range("M8").Value2 = 1.060884891993208
$ws.Range("N8").Value2 = 1.045733246780084

# Row 9
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.036967362378973
$ws.Range("D9").Value2 = 1.044631282742129
$ws.Range("E9").Value2 = 1.045639493486098
$ws.Range("F9").Value2 = 1.056023621706351
$ws.Range("I9").Value2 = 1.035325874907272
$ws.Range("J9").Value2 = 1.042760620361176
$ws.Range("K9").Value2 = 1.047765958877243
$ws.Range("L9").Value2 = 1.048770943347183
$ws.Range("M9").Value2 = 1.059122234220325
$ws.Range("N9").Value2 = 1.044241460048404

# Row 10
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.035431390097419
$ws.Range("D10").Value2 = 1.04343708094505
$ws.Range("E10").Value2 = 1.044287972053659
$ws.Range("F10").Value2 = 1.054567654148566
$ws.Range("I10").Value2 = 1.035073946131218
$ws.Range("J10").Value2 = 1.041767213000594
$ws.Range("K10").Value2 = 1.046856316463798
$ws.Range("L10").Value2 = 1.047704217703067
$ws.Range("M10").Value2 = 1.057948198507832
$ws.Range("N10").Value2 = 1.043246641935424

# Row 11
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.034766966413533
$ws.Range("D11").Value2 = 1.042920322578555
$ws.Range("E11").Value2 = 1.043703777967774
$ws.Range("F11").Value2 = 1.053938242267599
$ws.Range("I11").Value2 = 1.03496330840948
$ws.Range("J11").Value2 = 1.041336997628168
$ws.Range("K11").Value2 = 1.046461959520907
$ws.Range("L11").Value2 = 1.047242555332062
$ws.Range("M11").Value2 = 1.057440095261193
$ws.Range("N11").Value2 = 1.042815815607823

# Row 12
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.03452027002392
$ws.Range("D12").Value2 = 1.04272842767217
$ws.Range("E12").Value2 = 1.043486937224449
$ws.Range("F12").Value2 = 1.053704606895936
$ws.Range("I12").Value2 = 1.034921980091536
$ws.Range("J12").Value2 = 1.041177187761428
$ws.Range("K12").Value2 = 1.046315407241881
$ws.Range("L12").Value2 = 1.047071109877869
$ws.Range("M12").Value2 = 1.05725140350021
$ws.Range("N12").Value2 = 1.042655778792742

# Row 13
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.034573182694838
$ws.Range("D13").Value2 = 1.042769587400675
$ws.Range("E13").Value2 = 1.043533443269381
$ws.Range("F13").Value2 = 1.053754715386999
$ws.Range("I13").Value2 = 1.03493085567159
$ws.Range("J13").Value2 = 1.04121146790204
$ws.Range("K13").Value2 = 1.046346846380917
$ws.Range("L13").Value2 = 1.047107883843049
$ws.Range("M13").Value2 = 1.05729187667803
$ws.Range("N13").Value2 = 1.042690107615086

# Row 14
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.034746572366348
$ws.Range("D14").Value2 = 1.042904459414494
$ws.Range("E14").Value2 = 1.043685850676584
$ws.Range("F14").Value2 = 1.05391892670084
$ws.Range("I14").Value2 = 1.034959896938928
$ws.Range("J14").Value2 = 1.041323787867982
$ws.Range("K14").Value2 = 1.046449846897669
$ws.Range("L14").Value2 = 1.047228382847311
$ws.Range("M14").Value2 = 1.057424497110973
$ws.Range("N14").Value2 = 1.042802587088262

# Row 15
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.034853416755182
$ws.Range("D15").Value2 = 1.042987565462432
$ws.Range("E15").Value2 = 1.043779774477675
$ws.Range("F15").Value2 = 1.054020123449086
$ws.Range("I15").Value2 = 1.03497775942203
$ws.Range("J15").Value2 = 1.041392990774451
$ws.Range("K15").Value2 = 1.046513299594482
$ws.Range("L15").Value2 = 1.047302631125693
$ws.Range("M15").Value2 = 1.057506214322662
$ws.Range("N15").Value2 = 1.042871888270796

# Row 16
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.035475499719319
$ws.Range("D16").Value2 = 1.043471383742267
$ws.Range("E16").Value2 = 1.044326764756925
$ws.Range("F16").Value2 = 1.054609447983217
$ws.Range("I16").Value2 = 1.035081256156271
$ws.Range("J16").Value2 = 1.04179576371048
$ws.Range("K16").Value2 = 1.046882478690482
$ws.Range("L16").Value2 = 1.047734861783185
$ws.Range("M16").Value2 = 1.057981925255877
$ws.Range("N16").Value2 = 1.043275233190593

# Row 17
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.035865893723285
$ws.Range("D17").Value2 = 1.04377496177689
$ws.Range("E17").Value2 = 1.044670152180366
$ws.Range("F17").Value2 = 1.054979392769639
$ws.Range("I17").Value2 = 1.035145762030099
$ws.Range("J17").Value2 = 1.042048396119767
$ws.Range("K17").Value2 = 1.047113928260277
$ws.Range("L17").Value2 = 1.048006052720761
$ws.Range("M17").Value2 = 1.058280397032999
$ws.Range("N17").Value2 = 1.043528224366879

# Row 18
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.036093667815484
$ws.Range("D18").Value2 = 1.043952066323058
$ws.Range("E18").Value2 = 1.044870542940436
$ws.Range("F18").Value2 = 1.055195274764312
$ws.Range("I18").Value2 = 1.03518323754453
$ws.Range("J18").Value2 = 1.042195746155598
$ws.Range("K18").Value2 = 1.047248882910817
$ws.Range("L18").Value2 = 1.048164256599525
$ws.Range("M18").Value2 = 1.058454515776056
$ws.Range("N18").Value2 = 1.043675783656661

# Row 19
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.036171343710045
$ws.Range("D19").Value2 = 1.044012459890654
$ws.Range("E19").Value2 = 1.044938887714421
$ws.Range("F19").Value2 = 1.055268901778383
$ws.Range("I19").Value2 = 1.035195990325728
$ws.Range("J19").Value2 = 1.042245987638139
$ws.Range("K19").Value2 = 1.04729489113397
$ws.Range("L19").Value2 = 1.048218203872003
$ws.Range("M19").Value2 = 1.058513890011096
$ws.Range("N19").Value2 = 1.043726096487871

# Row 20
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.035824001529892
$ws.Range("D20").Value2 = 1.043742387344296
$ws.Range("E20").Value2 = 1.044633299743653
$ws.Range("F20").Value2 = 1.054939690898694
$ws.Range("I20").Value2 = 1.035138856638073
$ws.Range("J20").Value2 = 1.042021291685913
$ws.Range("K20").Value2 = 1.047089100659777
$ws.Range("L20").Value2 = 1.047976954135426
$ws.Range("M20").Value2 = 1.058248371221308
$ws.Range("N20").Value2 = 1.043501081441619

# Row 21
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.034695510660037
$ws.Range("D21").Value2 = 1.042864741518164
$ws.Range("E21").Value2 = 1.043640966207039
$ws.Range("F21").Value2 = 1.053870566239321
$ws.Range("I21").Value2 = 1.034951351425102
$ws.Range("J21").Value2 = 1.041290712671793
$ws.Range("K21").Value2 = 1.046419517752819
$ws.Range("L21").Value2 = 1.04719289787879
$ws.Range("M21").Value2 = 1.05738544257504
$ws.Range("N21").Value2 = 1.042769464921501

# Row 22
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.033986563411151
$ws.Range("D22").Value2 = 1.042313233149078
$ws.Range("E22").Value2 = 1.043017943805545
$ws.Range("F22").Value2 = 1.053199269613719
$ws.Range("I22").Value2 = 1.034832114269901
$ws.Range("J22").Value2 = 1.040831318586783
$ws.Range("K22").Value2 = 1.045998116673392
$ws.Range("L22").Value2 = 1.046700142258775
$ws.Range("M22").Value2 = 1.056843119395739
$ws.Range("N22").Value2 = 1.042309418444198

# Row 23
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.034362334597838
$ws.Range("D23").Value2 = 1.042605569091199
$ws.Range("E23").Value2 = 1.043348134346576
$ws.Range("F23").Value2 = 1.053555050511741
$ws.Range("I23").Value2 = 1.034895451494989
$ws.Range("J23").Value2 = 1.041074856558505
$ws.Range("K23").Value2 = 1.046221547696714
$ws.Range("L23").Value2 = 1.0469613409655
$ws.Range("M23").Value2 = 1.057130592613021
$ws.Range("N23").Value2 = 1.042553302267773

# Row 24
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.035842930611487
$ws.Range("D24").Value2 = 1.043757106227571
$ws.Range("E24").Value2 = 1.044649951468684
$ws.Range("F24").Value2 = 1.054957630158011
$ws.Range("I24").Value2 = 1.035141977349981
$ws.Range("J24").Value2 = 1.042033539032001
$ws.Range("K24").Value2 = 1.047100319325851
$ws.Range("L24").Value2 = 1.047990102463244
$ws.Range("M24").Value2 = 1.058262842229275
$ws.Range("N24").Value2 = 1.043513346180344

# Row 25
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.03756358325987
$ws.Range("D25").Value2 = 1.045094675847541
$ws.Range("E25").Value2 = 1.046164504772883
$ws.Range("F25").Value2 = 1.056589142903641
$ws.Range("I25").Value2 = 1.035422183920451
$ws.Range("J25").Value2 = 1.043145791525238
$ws.Range("K25").Value2 = 1.048118278277778
$ws.Range("L25").Value2 = 1.049184813693903
$ws.Range("M25").Value2 = 1.059577739803832
$ws.Range("N25").Value2 = 1.044627178199698
